$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "Mes" column (G), shifting
# Mes/Observaciones from G/H to I/J, and leaving room for the two new
# date columns.
$ws.Range("G1:H1").EntireColumn.Insert()

# New header cells
$ws.Range("G1").Value = "FechaRadicacion"
$ws.Range("H1").Value = "FechaMovimiento"

# Apply the custom date/time number format to the new date columns
# before assigning values, so Excel does not silently create an
# intermediate "m/d/yyyy" number format.
$dateFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G2:H5").NumberFormat = $dateFormat

function Set-DateCell($addr, [int]$year, [int]$month, [int]$day) {
    $d = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $ws.Range($addr).Value = $d
}

# Row 2
Set-DateCell "G2" 2025 8 1
Set-DateCell "H2" 2025 8 5
$ws.Range("I2").Value = "Agosto"
$ws.Range("J2").Value = "Falta soporte"

# Row 3
Set-DateCell "G3" 2025 8 2
Set-DateCell "H3" 2025 8 6
$ws.Range("I3").Value = "Agosto"

# Row 4
Set-DateCell "G4" 2025 9 3
Set-DateCell "H4" 2025 9 6
$ws.Range("I4").Value = "Septiembre"
$ws.Range("J4").Value = "Listo para radicar"

# Row 5
Set-DateCell "G5" 2025 9 4
Set-DateCell "H5" 2025 9 10
$ws.Range("I5").Value = "Septiembre"
$ws.Range("J5").Value = "Entregado"
